# chore: add monthly employment outputs
#
# 1) Refresh the "collected_at" timestamp (column J, rows 2-50) on every
#    sheet from 2026-02-12T23:04:40 -> 2026-02-12T23:15:00.
# 2) On the "피보험자수" sheet (sheet 6), rows 2-13 got refreshed
#    current_value numbers (column E) from a newer collection run, and a
#    few of the derived signals (F/G/I) moved as a result.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2026-02-12T23:15:00"

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    for ($r = 2; $r -le 50; $r++) {
        $ws.Cells.Item($r, 10).Value = $newTimestamp
    }
}

$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("E2").Value = 1435
$ws6.Range("E3").Value = 2758
$ws6.Range("E4").Value = 816

$ws6.Range("E5").Value = 1054
$ws6.Range("F5").Value = "정상"
$ws6.Range("G5").Value = "정상"
$ws6.Range("I5").Value = "정상"

$ws6.Range("E6").Value = 701
$ws6.Range("E7").Value = 431
$ws6.Range("E8").Value = 783
$ws6.Range("E9").Value = 1144
$ws6.Range("E10").Value = 1683
$ws6.Range("E11").Value = 134
$ws6.Range("E12").Value = 1061

$ws6.Range("E13").Value = 1361
$ws6.Range("F13").Value = "정상"
$ws6.Range("G13").Value = "정상"
$ws6.Range("I13").Value = "주의"
